$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain numbers (single decimal point) and would be
# auto-converted to numeric cells by Excel. Force those specific cells to Text format
# first so the assigned string (including trailing zeros) is preserved verbatim,
# matching the original inline-string cell content.
$ws.Range("D2").Value = "27.345.00"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "1.826.33"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -3.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.02"
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3697"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07248"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8654"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.14"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "1.823.13"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.673"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.349"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07079"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.85"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008888"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.19"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").Value = "27.367.34"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.160"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "2.052.83"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.18"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("E28").Value = "  +6.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.295"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.91"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7658"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.505"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.845"
$ws.Range("E35").Value = "  -6.06%  "
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.123"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01959"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05273"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.892"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.126"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1681"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5068"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.661"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.57"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.21"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4726"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.670"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.820"
$ws.Range("E51").Value = "  -3.47%  "
